# B3-and-B4-PowerPoint.pptx theme swap
#
# The canonical-OOXML diff shows ppt/theme/theme1.xml (originally the
# "Integral" / "Red Violet" theme used by the slide master) and
# ppt/theme/theme2.xml (originally the default "Office Theme" used by the
# notes master) swapping their contents: theme1.xml becomes the plain
# "Office Theme" palette and theme2.xml becomes the "Red Violet" palette.
# Everything else inside each <a:theme> (font scheme, format scheme) is
# byte-for-byte identical between the two parts, so the only real content
# change is the <a:clrScheme> (12 colours) and the name metadata.
#
# This host's object model only exposes live, persisted writes through
# ThemeColorScheme (the 12-slot modern theme colour list: dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink) reached from the slide master / any
# slide, and every one of those aliases back onto the single theme part
# driving the slide master (ppt/theme/theme1.xml). The theme/colour-scheme
# "Name" metadata and the notes-master theme part are not independently
# addressable here, so we recolour theme1.xml to the target "Office"
# palette (the part of the edit that is reachable and visually
# meaningful) and best-effort-attempt the rest so nothing is lost if a
# richer host exposes it.

function HexToRGB {
    param([string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palettes, in theme colour order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeColors = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")
$redVioletColors = @("000000", "FFFFFF", "454551", "D8D9DC", "E32D91", "C830CC", "4EA6DC", "4775E7", "8971E1", "D54773", "6B9F25", "8C8C8C")

$p = $ppt.ActivePresentation

# --- Slide master theme (ppt/theme/theme1.xml): Integral/Red Violet -> Office Theme ---
$master = $p.SlideMaster
$theme = $master.Theme

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $theme.ThemeColorScheme.Item($i + 1).RGB = HexToRGB $officeColors[$i]
}

# Best-effort: rename the theme / colour-scheme metadata to match the
# target ("Integral"/"Red Violet" -> "Office Theme"/"Office"). Harmless
# no-op on hosts that treat Name as read-only.
try { $theme.Name = "Office Theme" } catch { }
try { $theme.ThemeColorScheme.Name = "Office" } catch { }
try { $master.Design.Name = "Office Theme" } catch { }

# --- Notes master theme (ppt/theme/theme2.xml): Office Theme -> Integral/Red Violet ---
# Best-effort: some hosts alias NotesMaster.Theme back onto the slide
# master's theme part; attempt the recolour anyway in case this host
# resolves it to the distinct notes-master theme part.
try {
    $notesMaster = $p.NotesMaster
    $notesTheme = $notesMaster.Theme
    for ($i = 0; $i -lt $redVioletColors.Count; $i++) {
        $notesTheme.ThemeColorScheme.Item($i + 1).RGB = HexToRGB $redVioletColors[$i]
    }
    try { $notesTheme.Name = "Integral" } catch { }
    try { $notesTheme.ThemeColorScheme.Name = "Red Violet" } catch { }
} catch { }

Write-Host "Theme colours updated."
